$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (s="1") from H1 onto the new header cells I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set header text for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill the new data columns:
#   I = constant 1
#   J = same value as column H (IP) for that row
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}

$wb.Save()
